$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.3464964993005633, 86.29678392075563, 3.082599426703578, 6.48142807727062, 1, 96.20730792403039),
    @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126),
    @(0.06328177979961902, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 1.328214039578707),
    @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 3.811642989160245),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.3464964993005633, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 3.21235251628842),
    @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 0, 4.266567469086664),
    @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
